# Adds a new weekly column (AA) for the week of 04_05_2021,
# carrying the death-by-age-group figures, plus a totals formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new week
$ws.Range("AA1").Value2 = "04_05_2021"

# Per age-group counts for the new week
$ws.Range("AA2").Value2  = 1     # 0-9
$ws.Range("AA3").Value2  = 0     # 10-19
$ws.Range("AA4").Value2  = 0     # 20-29
$ws.Range("AA5").Value2  = 7     # 30-39
$ws.Range("AA6").Value2  = 8     # 40-49
$ws.Range("AA7").Value2  = 62    # 50-59
$ws.Range("AA8").Value2  = 208   # 60-69
$ws.Range("AA9").Value2  = 654   # 70-79
$ws.Range("AA10").Value2 = 991   # 80-89
$ws.Range("AA11").Value2 = 560   # 90+

# Total row, same pattern as the rest of the table
$ws.Range("AA12").Formula = "=SUM(AA2:AA11)"

# Keep the view scrolled roughly where it was, shifted for the new column
try {
    $win = $excel.ActiveWindow
    $win.ScrollColumn = 18
} catch {
}
